$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric but must be stored as TEXT
# (matches the source data, which keeps prices as literal strings, e.g.
# "588.67", preserving exact formatting/precision instead of letting
# Excel auto-convert it to a floating point number).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "63.294.48"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "3.071.23"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue "D5" "588.67"
$ws.Range("E5").Value = "  -0.58%  "
Set-TextValue "D6" "152.05"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue "D8" "0.548"
$ws.Range("E8").Value = "  +3.45%  "
$ws.Range("D9").Value = "3.069.77"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("E10").Value = "  -3.73%  "
Set-TextValue "D11" "5.87"
$ws.Range("E11").Value = "  +0.57%  "
Set-TextValue "D12" "0.461"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  -2.21%  "
Set-TextValue "D14" "37.18"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.579.90"
$ws.Range("E16").Value = "  -2.46%  "
Set-TextValue "D17" "7.21"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "63.325.84"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "3.071.32"
$ws.Range("E19").Value = "  -2.41%  "
Set-TextValue "D20" "475.45"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  -1.92%  "
Set-TextValue "D23" "7.53"
$ws.Range("E23").Value = "  +0.48%  "
Set-TextValue "D24" "2.37"
$ws.Range("E24").Value = "  +2.28%  "
Set-TextValue "D25" "13.02"
Set-TextValue "D26" "81.44"
$ws.Range("E26").Value = "  +0.24%  "
Set-TextValue "D27" "0.998"
$ws.Range("E27").Value = "  -0.25%  "
Set-TextValue "D28" "9.87"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "7.29"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "2.68"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D31" "1.00"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -2.63%  "
Set-TextValue "D33" "0.115"
$ws.Range("E33").Value = "  +4.30%  "
Set-TextValue "D34" "27.27"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D37" "6.12"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D38" "3.34"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  -4.02%  "
Set-TextValue "D40" "9.34"
$ws.Range("E40").Value = "  +1.09%  "
Set-TextValue "D41" "50.30"
$ws.Range("E41").Value = "  -2.06%  "
Set-TextValue "D42" "443.23"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("E44").Value = "  -2.85%  "
Set-TextValue "D45" "40.03"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "2.807.11"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("E47").Value = "  +1.69%  "
Set-TextValue "D48" "130.96"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +0.02%  "
Set-TextValue "D50" "25.11"
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D51" "2.26"
$ws.Range("E51").Value = "  -0.02%  "
